# Update countries & provincias Spain
#
# Applies the daily data refresh to the "Pais" worksheet:
#  - updates the "Datos actualizados..." timestamp cell
#  - refreshes case counters for several countries
#  - re-ranks three groups of countries whose case counts changed enough
#    to change their relative order in the (descending, by total cases)
#    list: Bulgaria / Republica de Macedonia,
#           Senegal / Letonia / Republica de Chipre,
#           Sri Lanka / Uruguay

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp header -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 29 de Abril de 2020 a las 17:22"

# --- Helper: write a full data row -------------------------------------
# Columns: A Pais, B Casos totales, C Nuevos casos, D Casos activos,
#          E Recuperados, F Casos criticos, G Muertes hoy, H Muertes
function Set-Row($Row, $Pais, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# --- Plain data refresh (same country, updated counters) --------------
Set-Row 4   "Estados Unidos"         1038490 2725 143117 835935 19098 172 59438
Set-Row 20  "Suiza"                  29407   143  22600  5104   185   4   1703
Set-Row 28  "Singapur"               15641   690  1188   14439  22    0   14
Set-Row 44  "Noruega"                7680    20   32     7441   44    1   207
Set-Row 74  "Azerbaiyan"             1766    49   1267   476    18    1   23

# --- Re-ranked groups (country order changes because case counts moved) --
# Bulgaria overtakes Republica de Macedonia
Set-Row 80  "Bulgaria"                1447 48 243 1140 38 6 64
Set-Row 81  "Republica de Macedonia"  1442 21 627 742  13 2 73

# Senegal overtakes Letonia and Republica de Chipre
Set-Row 92  "Senegal"                 882 59 315 558 1  0 9
Set-Row 93  "Letonia"                 849 13 348 486 4  2 15
Set-Row 94  "Republica de Chipre"     843 6  148 680 15 0 15

# Sri Lanka overtakes Uruguay
Set-Row 104 "Sri Lanka"               627 8 134 486 2  0 7
Set-Row 105 "Uruguay"                 625 5 394 216 11 0 15
